$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Avverkningsanmälningar")

# Update column C (row 2 through 185) from 45180 to 45181 (the "Förändrad" / changed date)
$ws.Range("C2:C185").Value = 45181
